# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Numeric-looking values are written with a leading apostrophe so Excel
# keeps them as text (matching the original inline-string cell content)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.186.24'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '3.533.36'
$ws.Range("E3").Value = '  +3.27%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''598.07'
$ws.Range("E5").Value = '  +2.01%  '

$ws.Range("D6").Value = '''138.46'
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").Value = '3.531.74'
$ws.Range("E7").Value = '  +3.23%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.492'

$ws.Range("E10").Value = '  +3.18%  '

$ws.Range("D11").Value = '''6.80'
$ws.Range("E11").Value = '  -6.41%  '

$ws.Range("D12").Value = '''0.387'
$ws.Range("E12").Value = '  +3.06%  '

$ws.Range("D13").Value = '4.135.36'
$ws.Range("E13").Value = '  +3.09%  '

$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").Value = '''27.16'
$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").Value = '3.537.51'
$ws.Range("E16").Value = '  +3.16%  '

$ws.Range("E17").Value = '  +1.46%  '

$ws.Range("D18").Value = '65.276.80'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("D19").Value = '''10.32'
$ws.Range("E19").Value = '  +4.91%  '

$ws.Range("D20").Value = '''5.95'
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").Value = '''14.30'
$ws.Range("E21").Value = '  +5.00%  '

$ws.Range("D22").Value = '''392.89'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").Value = '''0.573'
$ws.Range("E23").Value = '  +3.40%  '

$ws.Range("D24").Value = '3.676.14'
$ws.Range("E24").Value = '  +3.15%  '

$ws.Range("D25").Value = '''73.84'
$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("E27").Value = '  +8.51%  '

$ws.Range("D28").Value = '''7.77'
$ws.Range("E28").Value = '  +8.83%  '

$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").Value = '''2.30'
$ws.Range("E30").Value = '  +3.30%  '

$ws.Range("D31").Value = '''8.15'
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").Value = '3.544.77'
$ws.Range("E32").Value = '  +3.28%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("D34").Value = '''23.76'
$ws.Range("E34").Value = '  +3.45%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").Value = '''1.27'
$ws.Range("E36").Value = '  +9.22%  '

$ws.Range("D37").Value = '''6.97'
$ws.Range("E37").Value = '  +2.07%  '

$ws.Range("D38").Value = '''168.44'
$ws.Range("E38").Value = '  -2.22%  '

$ws.Range("E39").Value = '  +4.85%  '

$ws.Range("D40").Value = '''5.00'
$ws.Range("E40").Value = '  +5.07%  '

$ws.Range("D41").Value = '''0.0801'
$ws.Range("E41").Value = '  +4.90%  '

$ws.Range("D42").Value = '''0.824'
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("D43").Value = '''26.02'
$ws.Range("E43").Value = '  +16.17%  '

$ws.Range("D44").Value = '''42.75'
$ws.Range("E44").Value = '  -1.99%  '

$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("E46").Value = '  +0.43%  '

$ws.Range("E47").Value = '  +4.37%  '

$ws.Range("E48").Value = '  +6.57%  '

$ws.Range("D49").Value = '''6.79'
$ws.Range("E49").Value = '  +4.21%  '

$ws.Range("D50").Value = '2.389.98'
$ws.Range("E50").Value = '  +8.83%  '

$ws.Range("D51").Value = '''303.02'
$ws.Range("E51").Value = '  +7.11%  '
